$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set first day value (B4) to 2
$ws.Range("B4").Value = 2

# Update the active selection to J8 (matches the saved view state)
$ws.Range("J8").Select()
